$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.5415580391554872
$ws.Range("D2").Value = 0.5916544251648128

$ws.Range("C3").Value = 1.185175080187879
$ws.Range("D3").Value = 0.2441661456511059

$ws.Range("C4").Value = -0.08005982729591334
$ws.Range("D4").Value = 0.9366590635718302

$ws.Range("C5").Value = -1.224273436643161
$ws.Range("D5").Value = 0.2292648197711251

$ws.Range("C6").Value = 0.6113794900705095
$ws.Range("D6").Value = 0.5450148000051098

$ws.Range("C7").Value = -0.4625421506586854
$ws.Range("D7").Value = 0.6466392353969383

$ws.Range("C8").Value = -1.975695895920662
$ws.Range("D8").Value = 0.05635058087269962

$ws.Range("C9").Value = -0.8304602642932233
$ws.Range("D9").Value = 0.412072297270528

$ws.Range("C10").Value = -2.315261378112026
$ws.Range("D10").Value = 0.02676574212528093

$ws.Range("C11").Value = -1.114348540160877
$ws.Range("D11").Value = 0.2729452061442326
